# Add a "Model Group" column to the commercial-building reference export.
# The new column belongs to the COMDAT group and is inserted immediately
# before the existing "Township" column (which is part of LEGDAT), so
# Township shifts one column to the right (BK -> BL) and the new column
# takes its old slot (BK).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at BK, shifting the old BK column (Township) to BL.
$ws.Columns("BK:BK").Insert()

# Row 1 holds the source-table group name; the new column comes from COMDAT.
$ws.Range("BK1").Value = "COMDAT"

# Row 2 holds the actual field/column label.
$ws.Range("BK2").Value = "Model Group"
